# "data till 30Jan 8AM"
# Adds the 29-Jan-2021 (column AI) daily collection figures to the
# January-2021 sheet, adds a missing branch/area name on row 79, and
# updates the active selection to the last-entered cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New collection figures for 29-Jan-2021 (column AI) ---------------
# Most of these cells simply take the worksheet/row default style (s=3),
# which is what Excel applies automatically when a previously-empty cell
# gets a value typed into it.
$newValues = @{
    "AI5"  = 1500
    "AI10" = 7000
    "AI15" = 2000
    "AI16" = 800
    "AI23" = 3000
    "AI25" = 3000
    "AI31" = 1000
    "AI41" = 4500
    "AI44" = 3000
    "AI54" = 2000
    "AI55" = 5000
    "AI62" = 3000
    "AI65" = 10000
    "AI66" = 6000
    "AI67" = 4000
    "AI68" = 800
    "AI90" = 1000
    "AI91" = 5000
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# A handful of rows use a highlighted cell style (the same style already
# used on other data-entry cells in that row) instead of the plain
# default style, so copy the formatting from a sibling cell first and
# then fill in the value.
$styledValues = @{
    "AI6"  = @{ Value = 3000; CopyFrom = "H6"  }
    "AI73" = @{ Value = 3000; CopyFrom = "L73" }
    "AI82" = @{ Value = 5000; CopyFrom = "AC82" }
}

foreach ($addr in $styledValues.Keys) {
    $info = $styledValues[$addr]
    $ws.Range($info.CopyFrom).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $info.Value
}

$excel.CutCopyMode = $false

# --- Missing branch/area name on row 79 --------------------------------
$ws.Range("D79").Value = "BADRABAD"

# --- Update selection to the latest data-entry cell ---------------------
$ws.Activate()
$ws.Range("AI67").Select()
